# Preparations for downloading the CERNOX sensor calibration from Epics (cont).
#
# Drop the "M" suffix from the sensor names stored in column A, rows 66-75
# of "Feuil1" (e.g. TT683M -> TT683, TT690M -> TT690, ...), and move the
# current selection/view down to the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$renames = @{
    "A66" = "TT683"
    "A67" = "TT684"
    "A68" = "TT685"
    "A69" = "TT686"
    "A70" = "TT687"
    "A71" = "TT688"
    "A72" = "TT690"
    "A73" = "TT691"
    "A74" = "TT692"
    "A75" = "TT693"
}

foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}

# Update the view to match where the edits happened: scroll so row 49 is at
# the top and select A76 (the cell right after the edited block).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("A76").Select()
